$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: row number, year (text), count (number)
$rows = @(
    @(3,  "2024", 276),
    @(4,  "2023", 317),
    @(5,  "2022", 298),
    @(6,  "2021", 258),
    @(7,  "2020", 208),
    @(8,  "2019", 143),
    @(9,  "2018", 126),
    @(10, "2017", 107),
    @(11, "2016", 86),
    @(12, "2015", 73),
    @(13, "2014", 40),
    @(14, "2013", 29),
    @(15, "2012", 31),
    @(16, "2011", 26),
    @(17, "2010", 20),
    @(18, "2009", 13),
    @(19, "2008", 8),
    @(20, "2007", 8),
    @(21, "2006", 2),
    @(22, "2005", 6),
    @(23, "2003", 6),
    @(24, "2002", 1),
    @(25, "2001", 1),
    @(26, "1998", 1),
    @(27, "1997", 2),
    @(28, "1996", 1),
    @(29, "1994", 1)
)

# Make sure the full year column is formatted as text so values like
# "2024" are stored as strings and not auto-converted to numbers.
$ws.Range("A3:A29").NumberFormat = "@"

foreach ($r in $rows) {
    $rowNum = $r[0]
    $year = $r[1]
    $count = $r[2]
    $ws.Cells.Item($rowNum, 1).Value = $year
    $ws.Cells.Item($rowNum, 2).Value = $count
}
